$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 8, shifting the current rows 8-10 down to rows 10-12
$ws.Rows("8:9").Insert()

# Fill in the two newly inserted rows (8 and 9) with "Principiantes" / "general" entries
$ws.Range("A8").Value = "Torneo FEG"
$ws.Range("B8").Value = "Principiantes"
$ws.Range("C8").Value = "general"
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = "Esborraz, Juan Cruz"
$ws.Range("F8").Value = 38

$ws.Range("A9").Value = "Torneo FEG"
$ws.Range("B9").Value = "Principiantes"
$ws.Range("C9").Value = "general"
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = "Núñez, Valentino"
$ws.Range("F9").Value = 41

# Append a new row 13 with the "Juveniles" / "caballeros" entry
$ws.Range("A13").Value = "Torneo FEG"
$ws.Range("B13").Value = "Juveniles"
$ws.Range("C13").Value = "caballeros"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "Liberatori, Augusto"
$ws.Range("F13").Value = 75
